$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: 192.168.1.113 -> 127.0.0.1 (keep existing text style)
$ws.Range("C2").Value = "127.0.0.1"

# E2: keep same value 192.168.0.24, but now styled like C2 (text format)
$ws.Range("E2").Value = "192.168.0.24"
$ws.Range("E2").NumberFormat = $ws.Range("C2").NumberFormat

# Update selection to E2
$ws.Range("E2").Select()
